$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to retain its original text formatting (the source
# data uses dotted/locale-style numeric strings, many of which would otherwise
# be auto-coerced to numbers by COM value assignment). Apply a text number
# format before writing, then restore the default style so the saved file has
# no stray style index, exactly like the untouched cells around it.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.999.59'
$ws.Range("E2").Value = '  -0.61%  '

$ws.Range("D3").Value = '1.743.77'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '247.08'
$ws.Range("E5").Value = '  +1.85%  '

$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").Value = '0.5049'
$ws.Range("E7").Value = '  -4.93%  '

$ws.Range("E8").Value = '  -1.73%  '

$ws.Range("D9").Value = '0.06189'
$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").Value = '1.749.12'
$ws.Range("E10").Value = '  +0.24%  '

$ws.Range("E11").Value = '  +0.84%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '15.16'
$ws.Range("E12").Value = '  -1.05%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.6521'
$ws.Range("E13").Value = '  -0.11%  '

$ws.Range("D14").Value = '4.686'
$ws.Range("E14").Value = '  +0.96%  '

$ws.Range("D15").Value = '77.47'
$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.10%  '

$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").Value = '26.016.55'
$ws.Range("E18").Value = '  -0.48%  '

$ws.Range("D19").Value = '11.91'
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Value = '0.000006894'
$ws.Range("E20").Value = '  +1.22%  '

$ws.Range("D21").Value = '1.971.94'
$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").Value = '4.469'
$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("D23").Value = '8.743'
$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("E24").Value = '  +2.03%  '

$ws.Range("D25").Value = '135.87'
$ws.Range("E25").Value = '  -2.77%  '

$ws.Range("D26").Value = '1.507'
$ws.Range("E26").Value = '  -0.65%  '

$ws.Range("D27").Value = '15.28'
$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").Value = '1.778'
$ws.Range("E28").Value = '  -0.37%  '

$ws.Range("D29").Value = '105.72'
$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("D30").Value = '3.873'
$ws.Range("E30").Value = '  +1.91%  '

$ws.Range("D31").Value = '0.08196'
$ws.Range("E31").Value = '  -3.31%  '

$ws.Range("D32").Value = '3.643'
$ws.Range("E32").Value = '  -1.05%  '

$ws.Range("D33").Value = '0.04661'
$ws.Range("E33").Value = '  +1.11%  '

$ws.Range("D34").Value = '2.657'
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").Value = '0.9955'
$ws.Range("E35").Value = '  -0.61%  '

$ws.Range("D36").Value = '0.6084'
$ws.Range("E36").Value = '  -2.93%  '

$ws.Range("D37").Value = '2.791'
$ws.Range("E37").Value = '  +3.42%  '

$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D39").Value = '1.925'
$ws.Range("E39").Value = '  -1.18%  '

$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("D41").Value = '100.60'
$ws.Range("E41").Value = '  +0.78%  '

$ws.Range("D42").Value = '0.3923'
$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").Value = '0.7683'
$ws.Range("E43").Value = '  +2.03%  '

$ws.Range("D44").Value = '4.998'
$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("D45").Value = '0.1158'
$ws.Range("E45").Value = '  +0.73%  '

$ws.Range("D46").Value = '6.322'
$ws.Range("E46").Value = '  -0.16%  '

$ws.Range("D47").Value = '55.71'
$ws.Range("E47").Value = '  +1.41%  '

$ws.Range("D48").Value = '0.05324'
$ws.Range("E48").Value = '  -0.13%  '

$ws.Range("D49").Value = '30.68'
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("D50").Value = '7.609'
$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("D51").Value = '0.3465'
$ws.Range("E51").Value = '  -0.66%  '

# Restore default (General) style on the Price column so previously-untouched
# cells and the newly written ones all end up without an explicit style index.
$priceRange.Style = "Normal"
